$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 9015.583000000001
$ws.Range("J125").Value = 11872.889
$ws.Range("L125").Value = 106856.001
$ws.Range("N125").Value = -111776.001
$ws.Range("H141").Value = 2527.2
$ws.Range("I141").Value = 1961.9048
$ws.Range("J141").Value = 5495
$ws.Range("K141").Value = 5885.7144
$ws.Range("L141").Value = 16485
$ws.Range("M141").Value = -705.7143999999998
$ws.Range("N141").Value = -26845

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26090.902
$ws.Range("I32").Value = 30288.943
$ws.Range("K32").Value = 30288.943
$ws.Range("M32").Value = -30001.943
$ws.Range("H61").Value = 11222.973
$ws.Range("I61").Value = 8029.393
$ws.Range("J61").Value = 21158.555
$ws.Range("K61").Value = 8029.393
$ws.Range("L61").Value = 21158.555
$ws.Range("M61").Value = -7817.393
$ws.Range("N61").Value = -21582.555
$ws.Range("H63").Value = 3497.8572
$ws.Range("I63").Value = 2897
$ws.Range("K63").Value = 2897
$ws.Range("M63").Value = -2211
$ws.Range("H66").Value = 3497.8572
$ws.Range("I66").Value = 2897
$ws.Range("K66").Value = 14485
$ws.Range("M66").Value = -11053
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").ClearContents()
$ws.Range("N118").Value = 0
$ws.Range("H136").Value = 11222.973
$ws.Range("I136").Value = 8029.393
$ws.Range("J136").Value = 21158.555
$ws.Range("K136").Value = 24088.179
$ws.Range("L136").Value = 63475.665
$ws.Range("M136").Value = -21538.179
$ws.Range("N136").Value = -68575.66500000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0
$ws.Range("H87").Value = 220000
$ws.Range("J87").Value = 220000
$ws.Range("L87").Value = 220000
$ws.Range("N87").Value = -222496
$ws.Range("H90").Value = 220000
$ws.Range("J90").Value = 220000
$ws.Range("L90").Value = 660000
$ws.Range("N90").Value = -672480
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 510000000
$ws.Range("I6").Value = 510000000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 510000000
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -509999887
$ws.Range("H7").Value = 158.375
$ws.Range("I7").Value = 48.857143
$ws.Range("J7").Value = 243.55556
$ws.Range("K7").Value = 48.857143
$ws.Range("L7").Value = 243.55556
$ws.Range("M7").Value = 64.14285699999999
$ws.Range("N7").Value = -469.55556
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").Value = 3000
$ws.Range("N94").Value = -3902
$ws.Range("H115").Value = 37000
$ws.Range("J115").Value = 37000
$ws.Range("L115").Value = 37000
$ws.Range("N115").Value = -39350
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3380.3635
$ws.Range("J64").Value = 3993.3333
$ws.Range("L64").Value = 11979.9999
$ws.Range("N64").Value = -12519.9999
$ws.Range("H67").Value = 3380.3635
$ws.Range("J67").Value = 3993.3333
$ws.Range("L67").Value = 11979.9999
$ws.Range("N67").Value = -13851.9999
$ws.Range("H98").Value = 417.03845
$ws.Range("I98").Value = 329.4737
$ws.Range("J98").Value = 654.7143
$ws.Range("K98").Value = 988.4211
$ws.Range("L98").Value = 1964.1429
$ws.Range("M98").Value = 509.5789
$ws.Range("N98").Value = -4960.1429
$ws.Range("H131").Value = 33762.07
$ws.Range("J131").Value = 56720.47
$ws.Range("L131").Value = 170161.41
$ws.Range("N131").Value = -180241.41

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3000
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = 0
$ws.Range("H61").Value = 51176
$ws.Range("I61").Value = 67234.664
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 67234.664
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -67032.664
$ws.Range("N61").Value = -3404
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
$ws.Range("H113").Value = 51176
$ws.Range("I113").Value = 67234.664
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 67234.664
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -65064.664
$ws.Range("N113").Value = -7340

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 97750
$ws.Range("J47").Value = 97750
$ws.Range("L47").Value = 97750
$ws.Range("N47").Value = -98894
$ws.Range("H49").Value = 13829.6
$ws.Range("I49").Value = 2900
$ws.Range("J49").Value = 16562
$ws.Range("K49").Value = 2900
$ws.Range("L49").Value = 16562
$ws.Range("M49").Value = -2670
$ws.Range("N49").Value = -17022
